$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H5").Value = 203.55556
$ws.Range("I5").Value = 236
$ws.Range("J5").Value = 90
$ws.Range("K5").Value = 236
$ws.Range("L5").Value = 90
$ws.Range("M5").Value = -121
$ws.Range("N5").Value = -320
$ws.Range("H96").Value = 1758.2222
$ws.Range("I96").Value = 5250
$ws.Range("J96").Value = 760.5714
$ws.Range("K96").Value = 15750
$ws.Range("L96").Value = 2281.7142
$ws.Range("M96").Value = -14377
$ws.Range("N96").Value = -5027.7142
$ws.Range("H98").Value = 1478.1765
$ws.Range("I98").Value = 1087.7858
$ws.Range("K98").Value = 1087.7858
$ws.Range("M98").Value = 410.2141999999999
$ws.Range("H122").Value = 1478.1765
$ws.Range("I122").Value = 1087.7858
$ws.Range("K122").Value = 3263.3574
$ws.Range("M122").Value = -813.3574000000003
$ws.Range("H123").Value = 26113.334
$ws.Range("J123").Value = 26113.334
$ws.Range("L123").Value = 26113.334
$ws.Range("N123").Value = -35913.334
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 1289.1538
$ws.Range("I45").Value = 1178.125
$ws.Range("J45").Value = 1466.8
$ws.Range("K45").Value = 1178.125
$ws.Range("L45").Value = 1466.8
$ws.Range("M45").Value = -801.125
$ws.Range("N45").Value = -2220.8
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 281.4
$ws.Range("I22").Value = 279.92856
$ws.Range("J22").Value = 302
$ws.Range("K22").Value = 279.92856
$ws.Range("L22").Value = 302
$ws.Range("M22").Value = -106.92856
$ws.Range("N22").Value = -648
$ws.Range("H27").Value = 75500
$ws.Range("J27").Value = 75500
$ws.Range("L27").Value = 75500
$ws.Range("N27").Value = -75884
$ws.Range("H101").Value = 58700
$ws.Range("I101").Value = 39800
$ws.Range("J101").Value = 65000
$ws.Range("K101").Value = 39800
$ws.Range("L101").Value = 65000
$ws.Range("M101").Value = -36555
$ws.Range("N101").Value = -71490
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 72.21429000000001
$ws.Range("I7").Value = 74.8
$ws.Range("K7").Value = 74.8
$ws.Range("M7").Value = 38.2
$ws.Range("H22").Value = 62500970
$ws.Range("I22").Value = 71429540
$ws.Range("K22").Value = 71429540
$ws.Range("M22").Value = -71429190
$ws.Range("H32").Value = 7330
$ws.Range("I32").Value = 2042.8572
$ws.Range("K32").Value = 2042.8572
$ws.Range("M32").Value = -1726.8572
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 7387.2144
$ws.Range("I2").Value = 30
$ws.Range("J2").Value = 10330.1
$ws.Range("K2").Value = 180
$ws.Range("L2").Value = 61980.60000000001
$ws.Range("M2").Value = -67
$ws.Range("N2").Value = -62206.60000000001
$ws.Range("H5").Value = 427.9
$ws.Range("I5").Value = 397.66666
$ws.Range("K5").Value = 1192.99998
$ws.Range("M5").Value = -1080.99998
$ws.Range("H23").Value = 1831764.4
$ws.Range("I23").Value = 3663207.5
$ws.Range("J23").Value = 321.14285
$ws.Range("K23").Value = 10989622.5
$ws.Range("L23").Value = 963.4285500000001
$ws.Range("M23").Value = -10989387.5
$ws.Range("N23").Value = -1433.42855
$ws.Range("H69").Value = 800
$ws.Range("I69").Value = 800
$ws.Range("J69").Value = 0
$ws.Range("K69").Value = 2400
$ws.Range("L69").Value = 0
$ws.Range("M69").ClearContents()
$ws.Range("N69").Value = -1589
$ws.Range("H72").Value = 800
$ws.Range("I72").Value = 800
$ws.Range("J72").Value = 0
$ws.Range("K72").Value = 7200
$ws.Range("L72").Value = 0
$ws.Range("M72").ClearContents()
$ws.Range("N72").Value = -3144
$ws.Range("H131").Value = 898.2542
$ws.Range("I131").Value = 369.9
$ws.Range("J131").Value = 1006.0816
$ws.Range("K131").Value = 1109.7
$ws.Range("L131").Value = 3018.2448
$ws.Range("M131").Value = 3930.3
$ws.Range("N131").Value = -13098.2448
$ws.Range("H135").Value = 427.9
$ws.Range("I135").Value = 397.66666
$ws.Range("K135").Value = 3578.99994
$ws.Range("M135").Value = -1043.99994
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 1576.5555
$ws.Range("I122").Value = 1461.5454
$ws.Range("J122").Value = 1655.625
$ws.Range("K122").Value = 4384.6362
$ws.Range("L122").Value = 4966.875
$ws.Range("M122").Value = -1934.6362
$ws.Range("N122").Value = -9866.875
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H17").Value = 2675.6667
$ws.Range("I17").Value = 900.5
$ws.Range("J17").Value = 3182.8572
$ws.Range("K17").Value = 900.5
$ws.Range("L17").Value = 3182.8572
$ws.Range("M17").Value = -730.5
$ws.Range("N17").Value = -3522.8572
$ws.Range("H22").Value = 928.8570999999999
$ws.Range("I22").Value = 1025
$ws.Range("J22").Value = 800.6667
$ws.Range("K22").Value = 1025
$ws.Range("L22").Value = 800.6667
$ws.Range("M22").Value = -730
$ws.Range("N22").Value = -1390.6667
$ws.Range("H25").Value = 5535
$ws.Range("J25").Value = 5535
$ws.Range("L25").Value = 5535
$ws.Range("N25").Value = -5995
$ws.Range("H26").Value = 10000
$ws.Range("J26").Value = 10000
$ws.Range("L26").Value = 10000
$ws.Range("N26").Value = -10590
$ws.Range("H27").Value = 928.8570999999999
$ws.Range("I27").Value = 1025
$ws.Range("J27").Value = 800.6667
$ws.Range("K27").Value = 1025
$ws.Range("L27").Value = 800.6667
$ws.Range("M27").Value = -918
$ws.Range("N27").Value = -1014.6667
$ws.Range("H31").Value = 6166.5
$ws.Range("I31").Value = 999
$ws.Range("J31").Value = 7200
$ws.Range("K31").Value = 999
$ws.Range("L31").Value = 7200
$ws.Range("M31").Value = -751
$ws.Range("N31").Value = -7696
$ws.Range("H34").Value = 15560
$ws.Range("I34").Value = 6000
$ws.Range("J34").Value = 17950
$ws.Range("K34").Value = 6000
$ws.Range("L34").Value = 17950
$ws.Range("M34").Value = -5828
$ws.Range("N34").Value = -18294
$ws.Range("H122").Value = 11002909
$ws.Range("I122").Value = 1518057.9
$ws.Range("J122").Value = 29414680
$ws.Range("K122").Value = 4554173.699999999
$ws.Range("L122").Value = 88244040
$ws.Range("M122").Value = -4551723.699999999
$ws.Range("N122").Value = -88248940
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H29").Value = 5080
$ws.Range("J29").Value = 5080
$ws.Range("L29").Value = 5080
$ws.Range("N29").Value = -5660
$ws.Range("H136").Value = 3486.4358
$ws.Range("I136").Value = 691.8276
$ws.Range("J136").Value = 11590.8
$ws.Range("K136").Value = 2075.4828
$ws.Range("L136").Value = 34772.39999999999
$ws.Range("M136").Value = 474.5172000000002
$ws.Range("N136").Value = -39872.39999999999
